# VyTrackLoginPage DDT: add a new "VyTrackQA2User" sheet with the
# credentials used by the data-driven test and the PASSED/FAILED result
# column written back after the run.

$wb = $excel.ActiveWorkbook
$employees = $wb.Worksheets.Item(1)

# New sheet goes right after "Employees" and becomes the active tab.
$qa = $wb.Worksheets.Add($null, $employees)
$qa.Name = "VyTrackQA2User"

# Header row
$qa.Range("A1").Value = "username"
$qa.Range("B1").Value = "password"
$qa.Range("C1").Value = "firstname"
$qa.Range("D1").Value = "lastname"
$qa.Range("E1").Value = "Result"

# Row 2 - user1
$qa.Range("A2").Value = "user1"
$qa.Range("B2").Value = "UserUser123"
$qa.Range("C2").Value = "John"
$qa.Range("D2").Value = "Doe"
$qa.Range("E2").Value = "PASSED"

# Row 3 - user2
$qa.Range("A3").Value = "user2"
$qa.Range("B3").Value = "UserUser123"
$qa.Range("C3").Value = "Bella"
$qa.Range("D3").Value = "Stamm"
$qa.Range("E3").Value = "PASSED"

# Row 4 - storemanager51
$qa.Range("A4").Value = "storemanager51"
$qa.Range("B4").Value = "UserUser123"
$qa.Range("C4").Value = "Edd"
$qa.Range("D4").Value = "Turner"
$qa.Range("E4").Value = "PASSED"

# Row 5 - storemanager52 (name/password entered before username, as in the
# original authoring session)
$qa.Range("B5").Value = "UserUser123"
$qa.Range("C5").Value = "Roma"
$qa.Range("D5").Value = "Medhurst"
$qa.Range("A5").Value = "storemanager52"
$qa.Range("E5").Value = "PASSED"

# Row 6 - storemanager101
$qa.Range("A6").Value = "storemanager101"
$qa.Range("B6").Value = "UserUser123"
$qa.Range("C6").Value = "John"
$qa.Range("D6").Value = "Doe"
$qa.Range("E6").Value = "PASSED"

# Row 7 - storemanager102
$qa.Range("A7").Value = "storemanager102"
$qa.Range("B7").Value = "UserUser123"
$qa.Range("C7").Value = "John"
$qa.Range("D7").Value = "Doe"
$qa.Range("E7").Value = "PASSED"

# Column widths sized to fit the longest values in columns A and B
$qa.Columns.Item(1).ColumnWidth = 13.83
$qa.Columns.Item(2).ColumnWidth = 11

# Make the new sheet the active one, zoomed in, with C6:D7 selected -
# matches the view state the workbook was saved with.
$qa.Activate()
$excel.ActiveWindow.Zoom = 218
$qa.Range("C6:D7").Select()

Write-Output "VyTrackQA2User sheet added"
